$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation so date-like strings aren't converted to
# date serial numbers, then restore the default "Normal" style so no
# extra style index is left behind on the new cells.
$ws.Range("A4:D4").NumberFormat = "@"

$ws.Range("A4").Value = "2025-10-24 03:49:11"
$ws.Range("B4").Value = "2025-10-23"
$ws.Range("C4").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice23102025.pdf"
$ws.Range("D4").Value = "/home/runner/work/rashtriyametal_downloader/rashtriyametal_downloader/data/RashtriyaMetal/PDFs/ListPrice23102025.pdf"

$ws.Range("A4:D4").Style = "Normal"
